# Update the "Qty executed upto date" (C) values and the dependent
# "Upto date Amount" / "Amount Since prev bill" (G/H) text totals for the
# Bill Summary sheet, saved at the end of the loop iteration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty executed upto date (plain numeric cells) ---------------------------
$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 21
$ws.Range("C10").Value = 12
$ws.Range("C11").Value = 47
$ws.Range("C12").Value = 93
$ws.Range("C13").Value = 73
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 80
$ws.Range("C16").Value = 79
$ws.Range("C17").Value = 83

# --- Upto date Amount / Amount Since prev bill ------------------------------
# These cells store their number as literal text (no live formula in the
# sheet), so write the text via a throwaway formula and immediately collapse
# it back down to a plain cached value with Copy + PasteSpecial(values) -
# this keeps the cell's stored type/format identical to the original file
# instead of leaving a formula or reformatting the cell as Number/Text.
function Set-TextAmount($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Formula = "=""$text"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextAmount "G9"  "5376.00"
Set-TextAmount "G10" "5664.00"
Set-TextAmount "G11" "31114.00"
Set-TextAmount "G13" "9928.00"
Set-TextAmount "G14" "138.00"

Set-TextAmount "G19" "52220.00"
Set-TextAmount "H19" "52220.00"
Set-TextAmount "G21" "52220.00"
Set-TextAmount "H21" "52220.00"

$excel.CutCopyMode = $false
